# Updates the cryptocurrency price/volume snapshot on Sheet1 (and fixes the
# ImmutableX / MXToken row order swap) to match the latest scrape.
# Numeric-looking "Price" strings (single-dot decimals) are written with a
# temporary Text ("@") number format so Excel keeps them as exact text
# instead of silently coercing them to floating point numbers; the style is
# then reset to "Normal" so no stray cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.168.41'
$ws.Range('E2').Value = '  +3.50%  '
$ws.Range('D3').Value = '1.576.51'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  -0.88%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.08'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  -0.88%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.48'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.27%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('D12').Value = '1.801.92'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = '1.569.69'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('E14').Value = '  -0.87%  '
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('D16').Value = '28.124.12'
$ws.Range('E16').Value = '  +3.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.64'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.19%  '
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.45'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.73%  '
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.32'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.106'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('E33').Value = '  -1.03%  '
$ws.Range('D34').Value = '1.415.90'
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('E36').Value = '  -4.93%  '
$ws.Range('E37').Value = '  -1.17%  '
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.51'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.76%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.541'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('E45').Value = '  +4.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.63%  '
$ws.Range('D47').Value = '1.714.12'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.09'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.36%  '
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('E51').Value = '  -1.68%  '
